# Applies the cryptos list refresh (Thu Sep 28 18:20:20 UTC 2023 GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $text) {
    $range = $ws.Range($ref)
    if ($text -match "^-?\d+(\.\d+)?$") {
        # Purely numeric-looking text (e.g. '0.999') would otherwise be
        # auto-coerced to a number by Excel; force text with a leading
        # apostrophe, then restore the cell's original (default) style so
        # no stray number-format/quote-prefix formatting is left behind.
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

Set-TextValue 'D2' '27.190.90'
Set-TextValue 'E2' '  +3.86%  '
Set-TextValue 'D3' '1.662.46'
Set-TextValue 'E3' '  +4.64%  '
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.15%  '
Set-TextValue 'D5' '215.65'
Set-TextValue 'E5' '  +1.92%  '
Set-TextValue 'D6' '0.508'
Set-TextValue 'E6' '  +1.23%  '
Set-TextValue 'E7' '  -0.23%  '
Set-TextValue 'D8' '0.249'
Set-TextValue 'E8' '  +2.34%  '
Set-TextValue 'E9' '  +1.84%  '
Set-TextValue 'D10' '19.63'
Set-TextValue 'E10' '  +4.30%  '
Set-TextValue 'D11' '0.0862'
Set-TextValue 'E11' '  +1.28%  '
Set-TextValue 'D12' '1.896.19'
Set-TextValue 'E12' '  +4.52%  '
Set-TextValue 'D13' '1.671.30'
Set-TextValue 'E13' '  +5.00%  '
Set-TextValue 'E14' '  +1.95%  '
Set-TextValue 'E15' '  +3.70%  '
Set-TextValue 'D16' '64.89'
Set-TextValue 'E16' '  +2.25%  '
Set-TextValue 'B17' 'WrappedBTC'
Set-TextValue 'C17' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D17' '27.179.30'
Set-TextValue 'E17' '  +3.85%  '
Set-TextValue 'B18' 'BitcoinCash'
Set-TextValue 'C18' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D18' '239.06'
Set-TextValue 'E18' '  +5.08%  '
Set-TextValue 'D19' '7.88'
Set-TextValue 'E19' '  +4.60%  '
Set-TextValue 'E20' '  +1.69%  '
Set-TextValue 'E21' '  -0.20%  '
Set-TextValue 'D22' '4.44'
Set-TextValue 'E22' '  +5.42%  '
Set-TextValue 'D23' '2.27'
Set-TextValue 'E23' '  +5.34%  '
Set-TextValue 'D24' '9.32'
Set-TextValue 'E24' '  +5.03%  '
Set-TextValue 'D25' '146.09'
Set-TextValue 'E25' '  +0.51%  '
Set-TextValue 'E26' '  -0.16%  '
Set-TextValue 'D27' '7.19'
Set-TextValue 'E27' '  +3.76%  '
Set-TextValue 'E28' '  +1.65%  '
Set-TextValue 'D29' '15.84'
Set-TextValue 'E29' '  +3.69%  '
Set-TextValue 'E30' '  +1.34%  '
Set-TextValue 'E31' '  +1.54%  '
Set-TextValue 'D32' '1.544.12'
Set-TextValue 'E32' '  +6.42%  '
Set-TextValue 'E33' '  +3.24%  '
Set-TextValue 'D34' '3.05'
Set-TextValue 'E34' '  +3.68%  '
Set-TextValue 'D35' '1.58'
Set-TextValue 'E35' '  +8.61%  '
Set-TextValue 'E36' '  -0.10%  '
Set-TextValue 'D37' '0.575'
Set-TextValue 'E37' '  +2.22%  '
Set-TextValue 'D38' '0.892'
Set-TextValue 'E38' '  +9.50%  '
Set-TextValue 'E39' '  +3.26%  '
Set-TextValue 'D40' '5.96'
Set-TextValue 'E40' '  +3.80%  '
Set-TextValue 'E41' '  -0.14%  '
Set-TextValue 'D42' '2.28'
Set-TextValue 'E42' '  +5.36%  '
Set-TextValue 'D43' '66.49'
Set-TextValue 'E43' '  +10.56%  '
Set-TextValue 'D44' '1.803.31'
Set-TextValue 'E44' '  +4.30%  '
Set-TextValue 'D45' '0.774'
Set-TextValue 'E45' '  +2.52%  '
Set-TextValue 'D46' '0.922'
Set-TextValue 'E46' '  -0.74%  '
Set-TextValue 'D47' '90.21'
Set-TextValue 'E47' '  +3.18%  '
Set-TextValue 'B48' 'RenderToken'
Set-TextValue 'C48' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D48' '1.54'
Set-TextValue 'E48' '  +4.59%  '
Set-TextValue 'B49' 'BabyDogeCoin'
Set-TextValue 'C49' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D49' '0.0₆0105'
Set-TextValue 'E49' '  +7.10%  '
Set-TextValue 'E50' '  +0.99%  '
Set-TextValue 'D51' '0.0979'
Set-TextValue 'E51' '  +3.95%  '
